$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2: classical-best-embed vs. classical-best-tfidf
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.08799999999999999
$ws.Range("D2").Value = 0.043
$ws.Range("E2").Value = 0.041
$ws.Range("F2").Value = 0.017
$ws.Range("H2").Value = 0.049
$ws.Range("I2").Value = 0.047
$ws.Range("J2").Value = 0.048

# Row 3: BERT-base vs. classical-best-tfidf (label unchanged)
$ws.Range("D3").Value = 0.148
$ws.Range("E3").Value = 0.147
$ws.Range("F3").Value = 0.105
$ws.Range("G3").Value = 0.131
$ws.Range("H3").Value = 0.138
$ws.Range("I3").Value = 0.129
$ws.Range("J3").Value = 0.131

# Row 4: BERT-base vs. classical-best-embed
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.029
$ws.Range("D4").Value = 0.105
$ws.Range("E4").Value = 0.106
$ws.Range("F4").Value = 0.08799999999999999
$ws.Range("G4").Value = 0.081
$ws.Range("H4").Value = 0.089
$ws.Range("I4").Value = 0.082
$ws.Range("J4").Value = 0.083

# Row 5: BERT-base-nli vs. classical-best-tfidf (label unchanged)
$ws.Range("B5").Value = 0.256
$ws.Range("C5").Value = 0.203
$ws.Range("D5").Value = 0.167
$ws.Range("F5").Value = 0.118
$ws.Range("G5").Value = 0.115
$ws.Range("H5").Value = 0.119
$ws.Range("I5").Value = 0.162
$ws.Range("J5").Value = 0.147

# Row 6: BERT-base-nli vs. classical-best-embed
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.256
$ws.Range("C6").Value = 0.115
$ws.Range("D6").Value = 0.124
$ws.Range("E6").Value = 0.119
$ws.Range("F6").Value = 0.101
$ws.Range("G6").Value = 0.065
$ws.Range("H6").Value = 0.07000000000000001
$ws.Range("I6").Value = 0.115
$ws.Range("J6").Value = 0.099

# Row 7: BERT-base-nli vs. BERT-base (label unchanged)
$ws.Range("B7").Value = 0.256
$ws.Range("C7").Value = 0.08599999999999999
$ws.Range("D7").Value = 0.019
$ws.Range("E7").Value = 0.013
$ws.Range("F7").Value = 0.013
$ws.Range("G7").Value = -0.016
$ws.Range("H7").Value = -0.019
$ws.Range("I7").Value = 0.033
$ws.Range("J7").Value = 0.016
